# Updated cryptos list on Wed Sep 18 03:21:58 UTC 2024 with GitHub Actions
# Refresh live crypto prices / 1h volume percentages, and fix the
# USDe / EthereumClassic row ordering (rows 36-37 swapped back).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, inline pattern: force Text number format before writing a
# value that could otherwise be auto-coerced into a number by Excel's
# COM type inference (e.g. '548.35' or '0.999'), then drop back to the
# default 'Normal' style so no stray formatting is left behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '60.559.74'
Set-TextValue "E2" '  +4.23%  '
Set-TextValue "D3" '2.337.35'
Set-TextValue "E3" '  +2.25%  '
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '548.35'
Set-TextValue "E5" '  +2.75%  '
Set-TextValue "D6" '131.88'
Set-TextValue "E6" '  +0.72%  '
Set-TextValue "E7" '  -0.01%  '
Set-TextValue "D8" '0.582'
Set-TextValue "E8" '  -0.66%  '
Set-TextValue "D9" '2.335.48'
Set-TextValue "E9" '  +2.16%  '
Set-TextValue "E10" '  +1.67%  '
Set-TextValue "E11" '  +0.82%  '
Set-TextValue "E12" '  +0.25%  '
Set-TextValue "D13" '0.338'
Set-TextValue "E13" '  +1.79%  '
Set-TextValue "D14" '23.86'
Set-TextValue "E14" '  +2.00%  '
Set-TextValue "D15" '2.754.54'
Set-TextValue "E15" '  +2.28%  '
Set-TextValue "D16" '60.502.25'
Set-TextValue "E16" '  +4.26%  '
Set-TextValue "E17" '  +1.15%  '
Set-TextValue "D18" '2.345.73'
Set-TextValue "E18" '  +2.41%  '
Set-TextValue "D19" '10.65'
Set-TextValue "E19" '  +1.66%  '
Set-TextValue "D20" '4.17'
Set-TextValue "E20" '  -0.13%  '
Set-TextValue "D21" '315.93'
Set-TextValue "E21" '  +0.89%  '
Set-TextValue "D22" '6.70'
Set-TextValue "E22" '  +4.33%  '
Set-TextValue "E23" '  -0.24%  '
Set-TextValue "D24" '64.27'
Set-TextValue "E24" '  +2.19%  '
Set-TextValue "E25" '  +1.47%  '
Set-TextValue "D26" '0.999'
Set-TextValue "E26" '  -0.14%  '
Set-TextValue "D27" '7.87'
Set-TextValue "E27" '  -1.04%  '
Set-TextValue "D28" '1.37'
Set-TextValue "E28" '  +8.59%  '
Set-TextValue "D29" '1.21'
Set-TextValue "E29" '  +13.82%  '
Set-TextValue "D30" '173.63'
Set-TextValue "E30" '  +1.76%  '
Set-TextValue "E31" '  +2.91%  '
Set-TextValue "D32" '0.0₃0738'
Set-TextValue "E32" '  +2.35%  '
Set-TextValue "D33" '5.97'
Set-TextValue "E33" '  +3.72%  '
Set-TextValue "E34" '  +11.86%  '
Set-TextValue "D35" '0.381'
Set-TextValue "E35" '  +1.08%  '
Set-TextValue "B36" 'EthereumClassic'
Set-TextValue "C36" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D36" '17.99'
Set-TextValue "E36" '  +1.00%  '
Set-TextValue "B37" 'USDe'
Set-TextValue "C37" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D37" '0.999'
Set-TextValue "E37" '  +0.00%  '
Set-TextValue "E38" '  -0.05%  '
Set-TextValue "D39" '4.11'
Set-TextValue "E39" '  +5.61%  '
Set-TextValue "D40" '326.08'
Set-TextValue "E40" '  +13.86%  '
Set-TextValue "D41" '38.14'
Set-TextValue "E41" '  -0.44%  '
Set-TextValue "E42" '  +3.17%  '
Set-TextValue "D43" '140.27'
Set-TextValue "E43" '  +0.20%  '
Set-TextValue "D44" '3.49'
Set-TextValue "E44" '  +1.64%  '
Set-TextValue "D45" '0.0946'
Set-TextValue "E45" '  -0.71%  '
Set-TextValue "D46" '19.41'
Set-TextValue "E46" '  +7.89%  '
Set-TextValue "D47" '0.0498'
Set-TextValue "E47" '  +1.06%  '
Set-TextValue "D48" '0.563'
Set-TextValue "E48" '  +1.95%  '
Set-TextValue "D49" '0.0₆0224'
Set-TextValue "E49" '  +21.06%  '
Set-TextValue "E50" '  +2.11%  '
Set-TextValue "E51" '  +0.71%  '
